# Added cave graphic to template and code
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoardEnums")

# Extend the graphicEnums named range to cover the new cave graphic rows
$wb.Names.Item("graphicEnums").RefersTo = "=BoardEnums!`$A`$2:`$A`$28"

# Use an already-styled cell as the format source so the new cells pick up
# the same cellXf (wrap text, left/center alignment) already used in column G.
$ws.Range("G12").Copy()
foreach ($row in 11..27) {
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
}

# New graphic enum values, entered in the original authoring order (rows
# 11-13, 16-21, 23-25, 27 were typed first, then 15 and 14 were filled in).
$ws.Range("A11").Value = "cave_rock "
$ws.Range("A12").Value = "cave_rockLight "
$ws.Range("A13").Value = "cave_rockDark "
$ws.Range("A16").Value = "cave_rockMossy "
$ws.Range("A17").Value = "cave_rockWatery "
$ws.Range("A18").Value = "cave_rockWaterPuddle "
$ws.Range("A19").Value = "cave_rockBloody "
$ws.Range("A20").Value = "cave_rockBloodPuddle "
$ws.Range("A21").Value = "cave_rockLava "
$ws.Range("A23").Value = "cave_sandstone "
$ws.Range("A24").Value = "cave_sandstoneMossy "
$ws.Range("A25").Value = "cave_sandstoneWatery "
$ws.Range("A27").Value = "cave_lava"
$ws.Range("A15").Value = "cave_rockHole "
$ws.Range("A14").Value = "cave_rockBoulder"

# Column A widens to fit the new, longer graphic names
$ws.Range("A1").ColumnWidth = 30.5

# Matches the selection left behind in the authored workbook
$ws.Range("E17").Select() | Out-Null
